$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read current data rows (2-12, columns A-D) into memory
$numRows = 11
$data = @()
for ($i = 0; $i -lt $numRows; $i++) {
    $r = 2 + $i
    $row = @(
        $ws.Cells.Item($r, 1).Value(),
        $ws.Cells.Item($r, 2).Value(),
        $ws.Cells.Item($r, 3).Value(),
        $ws.Cells.Item($r, 4).Value()
    )
    $data += ,$row
}

# Sort rows ascending by column A (time)
$sorted = $data | Sort-Object { $_[0] }

# Write sorted rows back
for ($i = 0; $i -lt $numRows; $i++) {
    $r = 2 + $i
    $row = $sorted[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
